$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header values for the grade-distribution block (row 3, cols R:U) ---
$ws.Range("R3").Value = 2
$ws.Range("S3").Value = 3
$ws.Range("T3").Value = 4
$ws.Range("U3").Value = 5

# --- New formulas counting how many students got each grade (rows 4:32) ---
$ws.Range("R4:R32").Formula = "=IF(R`$3=`$M4,1,0)"
$ws.Range("S4:S32").Formula = "=IF(S`$3=`$M4,1,0)"
$ws.Range("T4:T32").Formula = "=IF(T`$3=`$M4,1,0)"
$ws.Range("U4:U32").Formula = "=IF(U`$3=`$M4,1,0)"

# --- Conditional formatting (3-color scale) for the new block ---
$cf = $ws.Range("R4:U32").FormatConditions.AddColorScale(3)
$cf.ColorScaleCriteria.Item(1).Type = 0
$cf.ColorScaleCriteria.Item(1).FormatColor.Color = 7039851
$cf.ColorScaleCriteria.Item(2).Type = 4
$cf.ColorScaleCriteria.Item(2).Value = 50
$cf.ColorScaleCriteria.Item(2).FormatColor.Color = 8711167
$cf.ColorScaleCriteria.Item(3).Type = 1
$cf.ColorScaleCriteria.Item(3).FormatColor.Color = 8109179

# --- View: zoom to 85%, scroll frozen pane back to top (C4), select T11 ---
$excel.ActiveWindow.Zoom = 85
$ws.Range("T11").Select()
